$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the whole "Bare conductive - Asked" paragraph (it is being moved
#    further down the list and reworded, see step 3).
# ---------------------------------------------------------------------------
$dashChar = [char]0x2013

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -eq "Bare conductive - Asked`r") {
        $candidate.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 2. Merge the three runs that make up "The Pi Hut " + "-" (en dash) + " Asked"
#    into a single run, same text/formatting, just no longer split up.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    $target = "The Pi Hut " + $dashChar + " Asked`r"
    if ($candidate.Range.Text -eq $target) {
        $piHutIndex = $i
        $mergedText = "The Pi Hut " + $dashChar + " Asked"
        $candidate.Range.Find.Execute($mergedText, $false, $false, $false, $false, $false, $true, 1, $false, $mergedText, 2) | Out-Null
        break
    }
}

# ---------------------------------------------------------------------------
# 3. The "_GoBack" bookmark currently sits at the end of "The Pi Hut" paragraph
#    -- it needs to move to the end of the brand-new paragraph created below.
#    Delete it now; it is re-created (in the correct spot) as part of the XML
#    inserted in step 4.
# ---------------------------------------------------------------------------
try {
    $d.Bookmarks.Item("_GoBack").Delete()
} catch {
    # no pre-existing _GoBack bookmark -- nothing to move, carry on
}

# ---------------------------------------------------------------------------
# 4. Insert a brand-new paragraph right after "The Pi Hut - Asked" containing
#    the re-worded "Bare conductive" prize entry, with the second run carrying
#    the web-paste formatting (Arial / grey / small caps shading) seen
#    elsewhere in this document (e.g. the "Ciseco" entry).
# ---------------------------------------------------------------------------
$piHutParagraph = $d.Paragraphs.Item($piHutIndex)
$piHutParagraph.Range.InsertParagraphAfter() | Out-Null
$newParagraph = $d.Paragraphs.Item($piHutIndex + 1)

$newParagraphXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve">Bare conductive - </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                <w:color w:val="222222"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
              </w:rPr>
              <w:t>10 pens, 2 house kits and 2 card kits</w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$newParagraph.Range.InsertXML($newParagraphXml) | Out-Null
